$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph whose text is exactly "Pose" (the last bullet
# under the "Character" heading). We search for the word "Pose" and
# collapse the found range to its end, i.e. right after the "P-o-s-e"
# characters but before the paragraph mark.
# ------------------------------------------------------------------
$poseRange = $d.Content
$found = $poseRange.Find.Execute("Pose", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Output "ERROR: could not find 'Pose' paragraph"
}
$poseRange.Collapse(0)   # wdCollapseEnd -> now a zero-length range right after "Pose"
$insertPos = $poseRange.Start

# ------------------------------------------------------------------
# Split the "Pose" paragraph into two paragraphs by inserting a new
# paragraph mark followed by the new bullet text "Move". A single
# trailing padding character ("Z") is appended after "Move" so that
# the later bookmark repositioning below lands on a range boundary
# that isn't the exact edge of this text insertion (repositioning a
# zero-width bookmark exactly onto either boundary of the just-typed
# text is unreliable) - the padding character is deleted again once
# the bookmark has been moved.
# ------------------------------------------------------------------
$splitRange = $d.Range($insertPos, $insertPos)
$splitRange.InsertAfter([char]13 + "MoveZ")

# Position right after "Move" (i.e. before the "Z" padding char / the
# new paragraph's paragraph mark) - this is where the _GoBack bookmark
# needs to end up, at the end of the newly created "Move" paragraph.
$newBookmarkPos = $insertPos + 1 + 4   # +1 for the paragraph mark, +4 for "Move"

# ------------------------------------------------------------------
# Move the "_GoBack" bookmark (originally right after "Pose") to the
# end of the new "Move" paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bookmarkTarget = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkTarget)

# Remove the single-character "Z" padding now that the bookmark has
# been anchored, leaving the new paragraph's text as plain "Move".
$padRange = $d.Range($newBookmarkPos, $newBookmarkPos + 1)
$padRange.Delete()

Write-Output "Inserted 'Move' bullet after 'Pose' and relocated the _GoBack bookmark."
